$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: add explicit <w:ind w:left="0" w:hanging="0"/> to the
# "Aquecendo na programação - Respostas" Titulo4 paragraph.
# -----------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Aquecendo na programa*") {
        $p.Format.LeftIndent = 0
        $p.Format.FirstLineIndent = -0.0001
    }
}

# -----------------------------------------------------------------
# Change 2: append three new paragraphs at the very end of the
# document (right before the final empty "Normal" paragraph's
# content, i.e. after it) :
#   1) an empty "Normal" paragraph
#   2) a "Normal" paragraph with "documentação do math.random"
#      (split into two runs, the second one carrying explicit
#      run formatting)
#   3) a "Normal" paragraph with the MDN url
# -----------------------------------------------------------------

# paragraph 1: plain empty paragraph, inherits formatting from the
# existing last paragraph (Normal / bidi=0 / jc=left / empty rPr).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# paragraph 2 (documentação do math.random) -- build it via a
# placeholder-character technique so the resulting paragraph keeps
# its correctly-inherited <w:pPr> (pStyle/bidi/jc/rPr) while still
# getting the exact run split + run-level formatting seen in the
# target document.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$pDocPara = $d.Paragraphs.Last
$pDocPara.Range.InsertAfter("X")

$insertPos = $pDocPara.Range.Start
$insertPoint = $d.Range($insertPos, $insertPos)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr/><w:t>documentaç</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Noto Serif CJK SC" w:cs="Lohit Devanagari"/><w:color w:val="auto"/><w:kern w:val="2"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>ão do math.random</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml)

# remove the temporary placeholder "X" character now found right
# before the paragraph mark.
$placeholderRange = $d.Range($pDocPara.Range.End - 2, $pDocPara.Range.End - 1)
$placeholderRange.Delete()

# paragraph 3: the MDN documentation url.
$pDocPara = $d.Paragraphs.Last
$pDocPara.Range.InsertParagraphAfter()
$pUrlPara = $d.Paragraphs.Last
$pUrlPara.Range.InsertAfter("https://developer.mozilla.org/pt-BR/docs/Web/JavaScript/Reference/Global_Objects/Math/random")

# -----------------------------------------------------------------
# Change 3: styles.xml Normal style overflowPunct false -> true
# -----------------------------------------------------------------
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.AutoAdjustRightIndent = $true
